$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

# Insert a new row at position 10, pushing the existing
# "master detail room" / "MST022" row down to row 11.
$ws.Rows.Item(10).Insert()

# New row 10: master room rate / MST022
$ws.Range("A10").Value = "master  room rate"
$ws.Range("B10").Value = "MST022"

# Row 11 (the shifted-down former row 10) gets its code updated to MST023
$ws.Range("A11").Value = "master detail room"
$ws.Range("B11").Value = "MST023"

$ws.Range("B12").Select() | Out-Null
